$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '28.227.40'
    'D3' = '1.800.59'
    'E3' = '  -0.90%  '
    'D4' = '1.004'
    'E4' = '  +0.20%  '
    'D5' = '314.26'
    'E5' = '  -0.32%  '
    'D6' = '1.003'
    'E6' = '  +0.18%  '
    'E7' = '  +3.30%  '
    'D8' = '0.3825'
    'E8' = '  -3.12%  '
    'D9' = '0.08002'
    'E9' = '  -0.94%  '
    'D10' = '41.30'
    'E10' = '  -0.95%  '
    'E11' = '  -0.61%  '
    'D12' = '6.313'
    'E12' = '  +0.44%  '
    'D13' = '1.004'
    'E13' = '  +0.25%  '
    'D14' = '20.60'
    'E14' = '  -1.81%  '
    'D15' = '1.805.93'
    'E15' = '  -0.58%  '
    'D16' = '7.321'
    'E16' = '  -2.36%  '
    'D17' = '92.03'
    'E17' = '  -0.51%  '
    'D18' = '0.00001091'
    'E18' = '  -3.90%  '
    'D19' = '0.06605'
    'E19' = '  -0.42%  '
    'D20' = '1.003'
    'D21' = '17.36'
    'E21' = '  -1.87%  '
    'D22' = '5.977'
    'E22' = '  -1.84%  '
    'D23' = '28.281.75'
    'E23' = '  -0.61%  '
    'D24' = '11.14'
    'E24' = '  -1.17%  '
    'D25' = '2.267'
    'E25' = '  -0.02%  '
    'D26' = '160.90'
    'E26' = '  +3.31%  '
    'D27' = '20.45'
    'E27' = '  -3.35%  '
    'D28' = '2.009.54'
    'E28' = '  -0.97%  '
    'D29' = '2.351'
    'E29' = '  -2.01%  '
    'D30' = '123.03'
    'E30' = '  -2.39%  '
    'D31' = '0.1084'
    'E31' = '  -1.33%  '
    'D32' = '1.058'
    'E32' = '  -4.33%  '
    'D33' = '3.687'
    'E33' = '  +0.97%  '
    'D34' = '5.552'
    'E34' = '  -4.11%  '
    'D35' = '0.07265'
    'E35' = '  +3.48%  '
    'D36' = '12.32'
    'E36' = '  +9.26%  '
    'D37' = '0.02309'
    'E37' = '  -0.86%  '
    'D38' = '0.2149'
    'E38' = '  -3.24%  '
    'D39' = '5.081'
    'E39' = '  -2.60%  '
    'D40' = '8.617'
    'E40' = '  -2.38%  '
    'D41' = '0.6193'
    'E41' = '  -1.46%  '
    'D42' = '1.166'
    'E42' = '  -0.71%  '
    'D43' = '1.370'
    'E43' = '  -2.24%  '
    'D44' = '13.28'
    'E44' = '  -1.36%  '
    'D45' = '0.6004'
    'E45' = '  +1.49%  '
    'E46' = '  +0.80%  '
    'D47' = '127.27'
    'E47' = '  +1.96%  '
    'D48' = '1.224'
    'E48' = '  +3.19%  '
    'D49' = '1.925'
    'E49' = '  -2.36%  '
    'D50' = '0.06802'
    'E50' = '  -1.22%  '
    'D51' = '73.02'
    'E51' = '  -1.47%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

Write-Host ("Updated " + $updates.Count + " cells")